# Update "想去人数" (F column) figures for the "展览" and "全部类型" sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Map: row number -> new value, for the "展览" sheet (sheet1)
$zhanlanUpdates = @{
    3  = 576
    4  = 1150
    6  = 103
    8  = 79
    9  = 1188
    10 = 16833
    11 = 301
    12 = 214
    14 = 6488
    15 = 658
    16 = 138
    18 = 36
    21 = 91
    23 = 644
    25 = 14
    26 = 24
    27 = 240
    28 = 917
    29 = 74
    30 = 5081
    33 = 11530
    34 = 1258
    35 = 27
    36 = 169
    37 = 233
    38 = 3868
    40 = 78
}

# Map: row number -> new value, for the "全部类型" sheet (sheet4)
$quanbuUpdates = @{
    3  = 576
    4  = 1150
    6  = 103
    8  = 79
    9  = 1188
    10 = 16833
    11 = 301
    12 = 214
    14 = 6488
    15 = 658
    16 = 138
    18 = 36
    21 = 91
    23 = 644
    25 = 14
    26 = 24
    27 = 240
    28 = 917
    29 = 74
    34 = 11530
    35 = 1258
    36 = 27
    37 = 169
    38 = 233
    39 = 3868
    41 = 78
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
